# EquipmentConf.xlsx maintenance edit.
#
# The sheet held equipment-config rows (id/name/quality/attack/defend) but was
# still carrying the default "Sheet1" tab name -- rename it to something that
# describes the data, and leave the cursor where the author last left it.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Name = "Properties"
$ws.Range("G25").Select()
